$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (was 45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update the prices in column D for rows 29-31
$ws.Range("D29").Value = 651
$ws.Range("D30").Value = 733
$ws.Range("D31").Value = 933
